$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = '大有能源'
$ws.Cells.Item(2, 2).Value = '大有能源'
$ws.Cells.Item(2, 3).Value = '大有能源'
$ws.Cells.Item(3, 1).Value = '科大国创'
$ws.Cells.Item(3, 2).Value = '黄河旋风'
$ws.Cells.Item(3, 3).Value = '科大国创'
$ws.Cells.Item(4, 1).Value = '黄河旋风'
$ws.Cells.Item(4, 2).Value = '山东墨龙'
$ws.Cells.Item(4, 3).Value = '黄河旋风'
$ws.Cells.Item(5, 1).Value = '神开股份'
$ws.Cells.Item(5, 2).Value = '中国核建'
$ws.Cells.Item(5, 3).Value = '平潭发展'
$ws.Cells.Item(6, 1).Value = '石化机械'
$ws.Cells.Item(6, 2).Value = '神开股份'
$ws.Cells.Item(6, 3).Value = '神开股份'
$ws.Cells.Item(7, 1).Value = '山东墨龙'
$ws.Cells.Item(7, 2).Value = '农业银行'
$ws.Cells.Item(7, 3).Value = '山东墨龙'
$ws.Cells.Item(8, 1).Value = '中国核建'
$ws.Cells.Item(8, 2).Value = '石化机械'
$ws.Cells.Item(8, 3).Value = '石化机械'
$ws.Cells.Item(9, 1).Value = '格尔软件'
$ws.Cells.Item(9, 2).Value = '郑州煤电'
$ws.Cells.Item(9, 3).Value = '多氟多'
$ws.Cells.Item(10, 1).Value = '多氟多'
$ws.Cells.Item(10, 2).Value = '闻泰科技'
$ws.Cells.Item(10, 3).Value = '神州信息'
$ws.Cells.Item(11, 1).Value = '神州信息'
$ws.Cells.Item(11, 2).Value = '盛新锂能'
$ws.Cells.Item(11, 3).Value = '赣锋锂业'
$ws.Cells.Item(12, 1).Value = '三花智控'
$ws.Cells.Item(12, 2).Value = '科大国创'
$ws.Cells.Item(12, 3).Value = '合肥城建'
$ws.Cells.Item(13, 1).Value = '赣锋锂业'
$ws.Cells.Item(13, 2).Value = '陕西黑猫'
$ws.Cells.Item(13, 3).Value = '三花智控'
$ws.Cells.Item(14, 1).Value = '闻泰科技'
$ws.Cells.Item(14, 2).Value = '三花智控'
$ws.Cells.Item(14, 3).Value = '中国核建'
$ws.Cells.Item(15, 1).Value = '盛新锂能'
$ws.Cells.Item(15, 2).Value = '多氟多'
$ws.Cells.Item(15, 3).Value = '安泰科技'
$ws.Cells.Item(16, 1).Value = '大众公用'
$ws.Cells.Item(16, 2).Value = '东方财富'
$ws.Cells.Item(16, 3).Value = '中信重工'
$ws.Cells.Item(17, 1).Value = '合肥城建'
$ws.Cells.Item(17, 2).Value = '赣锋锂业'
$ws.Cells.Item(17, 3).Value = '格尔软件'
$ws.Cells.Item(18, 1).Value = '农业银行'
$ws.Cells.Item(18, 2).Value = '中化岩土'
$ws.Cells.Item(18, 3).Value = '大众公用'
$ws.Cells.Item(19, 1).Value = '郑州煤电'
$ws.Cells.Item(19, 2).Value = '格尔软件'
$ws.Cells.Item(19, 3).Value = '华建集团'
$ws.Cells.Item(20, 1).Value = '国盾量子'
$ws.Cells.Item(20, 2).Value = '山西焦化'
$ws.Cells.Item(20, 3).Value = '闻泰科技'
$ws.Cells.Item(21, 1).Value = '平潭发展'
$ws.Cells.Item(21, 2).Value = '幸福蓝海'
$ws.Cells.Item(21, 3).Value = '盈新发展'
